# Update NATMI LR-pair stats for Vcan-Cd44 sheet (rows 2-26), following
# the re-analysis described in the commit "Natmi following Dr Hou advice".
# Ligand/receptor-expressing-cell counts (E, K) move from 1 to 3 replicate
# samples, and the average/total expression + specificity columns are
# refreshed to match the values recomputed from the new sample counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Cells = @{ "E"=3; "G"=8.877883333333335; "H"=26.63365; "I"=0.02762748111141948; "J"=0.02762748111141948; "K"=3; "M"=22.34478233333333; "N"=67.034347; "O"=0.03165884810812076; "P"=0.03165884810812076; "Q"=198.3743706640611; "R"=1785.36933597655; "S"=0.0008746542281164047; "T"=0.0008746542281164046 } }
    @{ Row = 3; Cells = @{ "E"=3; "G"=8.877883333333335; "H"=26.63365; "I"=0.02762748111141948; "J"=0.02762748111141948; "K"=3; "M"=92.44713066666667; "N"=277.341392; "O"=0.1309822411400946; "P"=0.1309822411400946; "Q"=820.7348405600891; "R"=7386.6135650408; "S"=0.003618709393029356; "T"=0.003618709393029355 } }
    @{ Row = 4; Cells = @{ "E"=3; "G"=8.877883333333335; "H"=26.63365; "I"=0.02762748111141948; "J"=0.02762748111141948; "K"=3; "M"=243.96462; "N"=731.89386; "O"=0.3456573768818275; "P"=0.3456573768818275; "Q"=2165.889433821; "R"=19493.004904389; "S"=0.009549642650825493; "T"=0.009549642650825493 } }
    @{ Row = 5; Cells = @{ "E"=3; "G"=8.877883333333335; "H"=26.63365; "I"=0.02762748111141948; "J"=0.02762748111141948; "K"=3; "M"=281.5837096666667; "N"=844.751129; "O"=0.3989573834764815; "P"=0.3989573834764815; "Q"=2499.867322987873; "R"=22498.80590689085; "S"=0.01102218757625783; "T"=0.01102218757625783 } }
    @{ Row = 6; Cells = @{ "E"=3; "G"=8.877883333333335; "H"=26.63365; "I"=0.02762748111141948; "J"=0.02762748111141948; "K"=3; "M"=65.45872566666667; "N"=196.376177; "O"=0.09274415039347572; "P"=0.09274415039347571; "Q"=581.134929617339; "R"=5230.214366556051; "S"=0.002562287263190398; "T"=0.002562287263190398 } }
    @{ Row = 7; Cells = @{ "E"=3; "G"=155.3909403333333; "H"=466.172821; "I"=0.4835679978836785; "J"=0.4835679978836785; "K"=3; "M"=22.34478233333333; "N"=67.034347; "O"=0.03165884810812076; "P"=0.03165884810812076; "Q"=3472.176738320321; "R"=31249.59064488289; "S"=0.01530920579494744; "T"=0.01530920579494744 } }
    @{ Row = 8; Cells = @{ "E"=3; "G"=155.3909403333333; "H"=466.172821; "I"=0.4835679978836785; "J"=0.4835679978836785; "K"=3; "M"=92.44713066666667; "N"=277.341392; "O"=0.1309822411400946; "P"=0.1309822411400946; "Q"=14365.44656541187; "R"=129289.0190887068; "S"=0.06333882010643274; "T"=0.06333882010643273 } }
    @{ Row = 9; Cells = @{ "E"=3; "G"=155.3909403333333; "H"=466.172821; "I"=0.4835679978836785; "J"=0.4835679978836785; "K"=3; "M"=243.96462; "N"=731.89386; "O"=0.3456573768818275; "P"=0.3456573768818275; "Q"=37909.89170986434; "R"=341189.0253887791; "S"=0.1671488456924694; "T"=0.1671488456924694 } }
    @{ Row = 10; Cells = @{ "E"=3; "G"=155.3909403333333; "H"=466.172821; "I"=0.4835679978836785; "J"=0.4835679978836785; "K"=3; "M"=281.5837096666667; "N"=844.751129; "O"=0.3989573834764815; "P"=0.3989573834764815; "Q"=43755.55742765166; "R"=393800.0168488649; "S"=0.1929230231686331; "T"=0.1929230231686331 } }
    @{ Row = 11; Cells = @{ "E"=3; "G"=155.3909403333333; "H"=466.172821; "I"=0.4835679978836785; "J"=0.4835679978836785; "K"=3; "M"=65.45872566666667; "N"=196.376177; "O"=0.09274415039347572; "P"=0.09274415039347571; "Q"=10171.69293436504; "R"=91545.23640928532; "S"=0.04484810312119582; "T"=0.04484810312119582 } }
    @{ Row = 12; Cells = @{ "E"=3; "G"=6.008960666666667; "H"=18.026882; "I"=0.0186995527069248; "J"=0.0186995527069248; "K"=3; "M"=22.34478233333333; "N"=67.034347; "O"=0.03165884810812076; "P"=0.03165884810812076; "Q"=134.2689181462282; "R"=1208.420263316054; "S"=0.0005920062988383308; "T"=0.0005920062988383308 } }
    @{ Row = 13; Cells = @{ "E"=3; "G"=6.008960666666667; "H"=18.026882; "I"=0.0186995527069248; "J"=0.0186995527069248; "K"=3; "M"=92.44713066666667; "N"=277.341392; "O"=0.1309822411400946; "P"=0.1309822411400946; "Q"=555.5111719221937; "R"=4999.600547299744; "S"=0.002449309321870334; "T"=0.002449309321870334 } }
    @{ Row = 14; Cells = @{ "E"=3; "G"=6.008960666666667; "H"=18.026882; "I"=0.0186995527069248; "J"=0.0186995527069248; "K"=3; "M"=243.96462; "N"=731.89386; "O"=0.3456573768818275; "P"=0.3456573768818275; "Q"=1465.97380563828; "R"=13193.76425074452; "S"=0.006463638337539104; "T"=0.006463638337539104 } }
    @{ Row = 15; Cells = @{ "E"=3; "G"=6.008960666666667; "H"=18.026882; "I"=0.0186995527069248; "J"=0.0186995527069248; "K"=3; "M"=281.5837096666667; "N"=844.751129; "O"=0.3989573834764815; "P"=0.3989573834764815; "Q"=1692.025435761087; "R"=15228.22892184978; "S"=0.007460324620135277; "T"=0.007460324620135276 } }
    @{ Row = 16; Cells = @{ "E"=3; "G"=6.008960666666667; "H"=18.026882; "I"=0.0186995527069248; "J"=0.0186995527069248; "K"=3; "M"=65.45872566666667; "N"=196.376177; "O"=0.09274415039347572; "P"=0.09274415039347571; "Q"=393.3389078211238; "R"=3540.050170390115; "S"=0.00173427412854176; "T"=0.00173427412854176 } }
    @{ Row = 17; Cells = @{ "E"=3; "G"=133.1709713333333; "H"=399.512914; "I"=0.4144206853098676; "J"=0.4144206853098676; "K"=3; "M"=22.34478233333333; "N"=67.034347; "O"=0.03165884810812076; "P"=0.03165884810812076; "Q"=2975.676367561906; "R"=26781.08730805716; "S"=0.01312008152908841; "T"=0.01312008152908841 } }
    @{ Row = 18; Cells = @{ "E"=3; "G"=133.1709713333333; "H"=399.512914; "I"=0.4144206853098676; "J"=0.4144206853098676; "K"=3; "M"=92.44713066666667; "N"=277.341392; "O"=0.1309822411400946; "P"=0.1309822411400946; "Q"=12311.27418785959; "R"=110801.4676907363; "S"=0.05428175013670035; "T"=0.05428175013670034 } }
    @{ Row = 19; Cells = @{ "E"=3; "G"=133.1709713333333; "H"=399.512914; "I"=0.4144206853098676; "J"=0.4144206853098676; "K"=3; "M"=243.96462; "N"=731.89386; "O"=0.3456573768818275; "P"=0.3456573768818275; "Q"=32489.00541636756; "R"=292401.0487473081; "S"=0.1432475670097781; "T"=0.1432475670097781 } }
    @{ Row = 20; Cells = @{ "E"=3; "G"=133.1709713333333; "H"=399.512914; "I"=0.4144206853098676; "J"=0.4144206853098676; "K"=3; "M"=281.5837096666667; "N"=844.751129; "O"=0.3989573834764815; "P"=0.3989573834764815; "Q"=37498.77612795332; "R"=337488.9851515799; "S"=0.1653361922697551; "T"=0.1653361922697551 } }
    @{ Row = 21; Cells = @{ "E"=3; "G"=133.1709713333333; "H"=399.512914; "I"=0.4144206853098676; "J"=0.4144206853098676; "K"=3; "M"=65.45872566666667; "N"=196.376177; "O"=0.09274415039347572; "P"=0.09274415039347571; "Q"=8717.202079272198; "R"=78454.81871344979; "S"=0.03843509436454563; "T"=0.03843509436454563 } }
    @{ Row = 22; Cells = @{ "E"=3; "G"=17.89372566666666; "H"=53.681177; "I"=0.05568428298810961; "J"=0.05568428298810962; "K"=3; "M"=22.34478233333333; "N"=67.034347; "O"=0.03165884810812076; "P"=0.03165884810812076; "Q"=399.8314051540465; "R"=3598.482646386419; "S"=0.001762900257130175; "T"=0.001762900257130175 } }
    @{ Row = 23; Cells = @{ "E"=3; "G"=17.89372566666666; "H"=53.681177; "I"=0.05568428298810961; "J"=0.05568428298810962; "K"=3; "M"=92.44713066666667; "N"=277.341392; "O"=0.1309822411400946; "P"=0.1309822411400946; "Q"=1654.22359481982; "R"=14888.01235337838; "S"=0.007293652182061842; "T"=0.007293652182061842 } }
    @{ Row = 24; Cells = @{ "E"=3; "G"=17.89372566666666; "H"=53.681177; "I"=0.05568428298810961; "J"=0.05568428298810962; "K"=3; "M"=243.96462; "N"=731.89386; "O"=0.3456573768818275; "P"=0.3456573768818275; "Q"=4365.43598265258; "R"=39288.92384387322; "S"=0.01924768319121534; "T"=0.01924768319121534 } }
    @{ Row = 25; Cells = @{ "E"=3; "G"=17.89372566666666; "H"=53.681177; "I"=0.05568428298810961; "J"=0.05568428298810962; "K"=3; "M"=281.5837096666667; "N"=844.751129; "O"=0.3989573834764815; "P"=0.3989573834764815; "Q"=5038.581652977648; "R"=45347.23487679883; "S"=0.02221565584170016; "T"=0.02221565584170016 } }
    @{ Row = 26; Cells = @{ "E"=3; "G"=17.89372566666666; "H"=53.681177; "I"=0.05568428298810961; "J"=0.05568428298810962; "K"=3; "M"=65.45872566666667; "N"=196.376177; "O"=0.09274415039347572; "P"=0.09274415039347571; "Q"=1171.300479568925; "R"=10541.70431612033; "S"=0.0051643915160021; "T"=0.0051643915160021 } }
)

foreach ($update in $updates) {
    $r = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $ws.Range("$col$r").Value = $update.Cells[$col]
    }
}
